# Smokeview script-error slide: swap the rounded-rectangle "card" for a
# plain square card, and blow up / re-flow the error text.
#
# PowerPoint's Shape.Left/Top/Width/Height are Single-precision (float32)
# properties expressed in points, while the OOXML stores EMU (914400 per
# inch) integers. A naive EMU -> pt conversion can therefore truncate to
# one EMU below the intended value once it round-trips back through the
# Single. EmuToPt nudges the point value by the smallest amount needed so
# it reproduces the exact target EMU once PowerPoint stores it as a float32.
function EmuToPt($emu) {
    $pt = $emu * 72.0 / 914400.0
    for ($i = 0; $i -lt 1000; $i++) {
        $f32 = [single]$pt
        $back = [int64]([double]$f32 * 914400.0 / 72.0)
        if ($back -eq $emu) {
            return $pt
        }
        $pt += 0.0000001
    }
    return $pt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1: "Rectangle: Rounded Corners 1" -> "Rectangle 1" -------------
$card = $s.Shapes.Item(1)

$card.Name = "Rectangle 1"

# roundRect -> rect (msoShapeRectangle = 1)
$card.AutoShapeType = 1

# New position / size
$card.Left   = EmuToPt 751114
$card.Top    = EmuToPt 457199
$card.Width  = EmuToPt 4572000
$card.Height = EmuToPt 4572000

# Drop the heavy outline the card used to have
$card.Line.Visible = $false

# --- Shape 2: "TextBox 3" --------------------------------------------------
$label = $s.Shapes.Item(2)

# Update text + font size *before* the frame geometry: the textbox has
# spAutoFit, so touching the text re-derives Height from the (old) font
# size if it happens after we set the final size.
$tr = $label.TextFrame.TextRange
$tr.Text = "viewpoint`rnot found"
$tr.Font.Size = 72

$label.Left   = EmuToPt 1149152
$label.Top    = EmuToPt 1589037
$label.Width  = EmuToPt 10069286
$label.Height = EmuToPt 2308324

Write-Host "script error slide updated"
